$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 61 with data for 2020-07-30
$ws.Range("A61").NumberFormat = "@"
$ws.Range("A61").Value = "2020-07-30"
$ws.Range("A61").Style = "Normal"
$ws.Range("B61").Value = 416179
$ws.Range("C61").Value = 461775
$ws.Range("D61").Value = 90582
$ws.Range("E61").Value = 46000
$ws.Range("F61").Value = 27.38
